$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at position 304, pushing the existing rows 304-341 down to 307-344
$ws.Range("A304:A306").EntireRow.Insert()

# Populate the 3 newly inserted rows with the new weekly data (week of 2023-09-11, serial 45180)
# Row 304: Especial, bandeja 10 kilos
$ws.Cells.Item(304, 1).Value = 9
$ws.Cells.Item(304, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(304, 3).Value = "Metropolitana"
$ws.Cells.Item(304, 4).Value = 45180
$ws.Cells.Item(304, 5).Value = 13
$ws.Cells.Item(304, 6).Value = "Fruta"
$ws.Cells.Item(304, 7).Value = 100107
$ws.Cells.Item(304, 8).Value = "Otros"
$ws.Cells.Item(304, 9).Value = 100107002
$ws.Cells.Item(304, 10).Value = "Chirimoya"
$ws.Cells.Item(304, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(304, 12).Value = "Especial"
$ws.Cells.Item(304, 13).Value = 200
$ws.Cells.Item(304, 14).Value = 26000
$ws.Cells.Item(304, 15).Value = 26000
$ws.Cells.Item(304, 16).Value = 26000
$ws.Cells.Item(304, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(304, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(304, 19).Value = 2600
$ws.Cells.Item(304, 20).Value = 10

# Row 305: Extra (doble especial), bandeja 10 kilos
$ws.Cells.Item(305, 1).Value = 9
$ws.Cells.Item(305, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(305, 3).Value = "Metropolitana"
$ws.Cells.Item(305, 4).Value = 45180
$ws.Cells.Item(305, 5).Value = 13
$ws.Cells.Item(305, 6).Value = "Fruta"
$ws.Cells.Item(305, 7).Value = 100107
$ws.Cells.Item(305, 8).Value = "Otros"
$ws.Cells.Item(305, 9).Value = 100107002
$ws.Cells.Item(305, 10).Value = "Chirimoya"
$ws.Cells.Item(305, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(305, 12).Value = "Extra (doble especial)"
$ws.Cells.Item(305, 13).Value = 180
$ws.Cells.Item(305, 14).Value = 28000
$ws.Cells.Item(305, 15).Value = 28000
$ws.Cells.Item(305, 16).Value = 28000
$ws.Cells.Item(305, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(305, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(305, 19).Value = 2800
$ws.Cells.Item(305, 20).Value = 10

# Row 306: Primera, bandeja 10 kilos
$ws.Cells.Item(306, 1).Value = 9
$ws.Cells.Item(306, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(306, 3).Value = "Metropolitana"
$ws.Cells.Item(306, 4).Value = 45180
$ws.Cells.Item(306, 5).Value = 13
$ws.Cells.Item(306, 6).Value = "Fruta"
$ws.Cells.Item(306, 7).Value = 100107
$ws.Cells.Item(306, 8).Value = "Otros"
$ws.Cells.Item(306, 9).Value = 100107002
$ws.Cells.Item(306, 10).Value = "Chirimoya"
$ws.Cells.Item(306, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(306, 12).Value = "Primera"
$ws.Cells.Item(306, 13).Value = 180
$ws.Cells.Item(306, 14).Value = 24000
$ws.Cells.Item(306, 15).Value = 24000
$ws.Cells.Item(306, 16).Value = 24000
$ws.Cells.Item(306, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(306, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(306, 19).Value = 2400
$ws.Cells.Item(306, 20).Value = 10
